$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Step 1: fix the "rel" bullet paragraph (paragraph 45):
#   "defines relationship between current docume" + bookmark + "nt and external resource"
#   -> "defines relationship between current document and external resource" (no bookmark)
# ------------------------------------------------------------------
$relPara = $d.Paragraphs.Item(45)
$relRange = $relPara.Range
$relRange.Find.Execute("current docume", $true, $false, $false, $false, $false, $true, 1, $false, "current document", 2)

$relPara = $d.Paragraphs.Item(45)
$relRange2 = $relPara.Range
$relRange2.Find.Execute("nt and external resource", $true, $false, $false, $false, $false, $true, 1, $false, " and external resource", 2)

# ------------------------------------------------------------------
# Helper: appends a new plain ListParagraph bullet after $afterPara
# with the given 0-based list level and text. Returns the new
# Word.Paragraph object. No character formatting is touched here so
# formatting never "leaks" forward from a previous bullet.
# ------------------------------------------------------------------
function Add-Bullet($afterPara, $ilvl, $text) {
    $afterPara.Range.InsertParagraphAfter()
    $newIndex = $afterPara.Index + 1
    $newPara = $d.Paragraphs.Item($newIndex)
    $newPara.Range.ListFormat.ListLevelNumber = $ilvl + 1
    $s = $newPara.Range.Start
    $ins = $d.Range($s, $s)
    $ins.InsertAfter($text)
    return $d.Paragraphs.Item($newIndex)
}

# ------------------------------------------------------------------
# Step 2: insert the new (plain-text) paragraphs after the "rel"
# bullet, in document order. Character formatting (italic/underline
# headers) is applied afterwards in a separate pass so that it never
# leaks onto the paragraph that gets typed next.
# ------------------------------------------------------------------
$cur = $d.Paragraphs.Item(45)

$textFormattingPara = Add-Bullet $cur 0 "Text formatting"
$cur = $textFormattingPara
$cur = Add-Bullet $cur 1 "<b></b> = defines contained text as bold"
$cur = Add-Bullet $cur 1 "<i></i> = defines contained text as italicized"
$cur = Add-Bullet $cur 1 "<u></u> = defines contained text as underlined"
$cur = Add-Bullet $cur 1 "<del></del> = defines contained text as deleted"
$cur = Add-Bullet $cur 1 "<big></big> = defines contained text as big"
$cur = Add-Bullet $cur 1 "<small></small> = defines contained text as small"
$cur = Add-Bullet $cur 1 "<sub></sub> = defines contained text as subscripted"
$cur = Add-Bullet $cur 1 "<sup></sup> = defines contained text as superscripted"
$cur = Add-Bullet $cur 1 "<tt></tt> = defines contained text as monospaced"
$cur = Add-Bullet $cur 2 "style = assigns a unique style to the element"
$spanDivPara = Add-Bullet $cur 0 "Span & Div"
$cur = $spanDivPara
$spanPara = Add-Bullet $cur 1 "<span></span> = inline container to group elements for styling purposes"
$cur = $spanPara
$cur = Add-Bullet $cur 1 "<div></div> = block container to group elements for styling purposes"

# ------------------------------------------------------------------
# Step 3: now that every paragraph exists, go back and italicize +
# underline the two section-header runs ("Text formatting" and
# "Span & Div"). Doing this last means it can no longer leak onto
# later bullets via the "continue previous formatting" behavior of
# InsertParagraphAfter/InsertAfter.
# ------------------------------------------------------------------
function Format-Header($para, $text) {
    $s = $para.Range.Start
    $e = $s + $text.Length
    $fr = $d.Range($s, $e)
    $fr.Italic = 1
    $fr.Underline = 1
}

Format-Header $textFormattingPara "Text formatting"
Format-Header $spanDivPara "Span & Div"

# ------------------------------------------------------------------
# Step 4: move the "_GoBack" bookmark so it sits between "</span>"
# and " = inline container..." in the span bullet, matching the
# location Word leaves it at after the final edit of the document.
# ------------------------------------------------------------------
$spanEnd = $spanPara.Range.Start + "<span></span>".Length
$bmRange = $d.Range($spanEnd, $spanEnd)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

Write-Output "done"
